$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scrape timestamp applied to every data row (2..398).
$newTimestamp = "2022-12-29 12:55:52"

for ($row = 2; $row -le 398; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}

# ratingAmount (column D) bumped by one for a handful of products that
# picked up an extra rating between scrapes.
$ws.Range("D19").Value = 35
$ws.Range("D23").Value = 68
$ws.Range("D30").Value = 37
$ws.Range("D41").Value = 114
$ws.Range("D209").Value = 66
